$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D..J to F..L)
$ws.Range("D1:E1").EntireColumn.Insert()

# Set header text and copy style from existing header cell (C1) for the new headers
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats

# New sum_SASA / max_SASA values per row (rows 2-24)
$sumSasa = @{
    2  = 4.65941418339056
    3  = 5.076088470882665
    4  = 4.502203803170425
    5  = 2.651420443376771
    6  = 2.351384937873273
    7  = 2.679408235125698
    8  = 6.338615601748595
    9  = 6.843748702262895
    10 = 7.238006874690403
    11 = 7.341690157828683
    12 = 6.590928259043968
    13 = 7.151842652483401
    14 = 6.692305731388845
    15 = 6.627845993509669
    16 = 7.31901792606446
    17 = 7.478598491147868
    18 = 5.32610747010928
    19 = 5.053632587327314
    20 = 5.224285506329593
    21 = 4.610734076466935
    22 = 5.020146705401829
    23 = 4.58744599634927
    24 = 4.622398868120259
}

$maxSasa = @{
    2  = 2.362150809882157
    3  = 2.753384233219247
    4  = 2.280227683002259
    5  = 2.651420443376771
    6  = 2.351384937873273
    7  = 2.679408235125698
    8  = 2.251887232513918
    9  = 2.488910808023066
    10 = 2.505235247011588
    11 = 2.568130278967355
    12 = 2.402306252462433
    13 = 2.456060678658576
    14 = 2.423617038821429
    15 = 2.525591632529191
    16 = 2.560270440720127
    17 = 2.575386834662778
    18 = 2.76768229887049
    19 = 2.530553721246874
    20 = 2.637590683923242
    21 = 2.415892868585273
    22 = 2.545281005075958
    23 = 2.381465013269206
    24 = 2.325370683812271
}

for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 4).Value = $sumSasa[$row]
    $ws.Cells.Item($row, 5).Value = $maxSasa[$row]
}
